$wb = $excel.ActiveWorkbook

# Build the updated Cypher query string for StatOutput_Message (Akita -> Mastiff)
$mastiffStatQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Mastiff']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# StatOutput sheet: the number_of_files / number_of_sample counts dropped to 0 for the Mastiff filter
$statOutput = $wb.Worksheets.Item("StatOutput")
$statRange = $statOutput.Range("A2:B2")
$statRange.NumberFormat = "@"
$statRange.Value = "0"

# StatOutput_Message sheet: update the logged Cypher text to reflect the Mastiff breed filter
$statOutputMessage = $wb.Worksheets.Item("StatOutput_Message")
$statOutputMessage.Range("A18").Value = $mastiffStatQuery
